# Apply new "dSF" (column F) values to rows 2-25, as re-pulled/recomputed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 1
    4  = 4
    5  = -5
    6  = 2
    7  = -1
    8  = 5
    9  = -5
    10 = 0
    11 = -3
    12 = -1
    13 = 1
    15 = -5
    16 = -5
    17 = -2
    18 = 3
    20 = -3
    21 = -1
    22 = -5
    23 = -7
    25 = -1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $newValues[$row]
}
